# Standardise "cost_variable" -> "cost_variable_om" across the sheet,
# then restore the user selection to reflect the edited range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Replace all occurrences of the old parameter name with the new,
# standardised one (column C, rows 10-39 in this sheet).
$found = $ws.Cells.Replace("cost_variable", "cost_variable_om", 1, 1, $false, $false, $false)

# Reflect the edited range as the active selection, matching the
# author's final selection after making the change.
$ws.Range("C10:C39").Select()
